$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above row 15 (shifts old rows 15-19 down to 16-20) ---
$ws.Rows("15").Insert()

# --- 2. Fill the new row 15 with the updated "delzescaux2023lagrange" record
#         (previously a "preprint" titled "Lagrange multiplier field approach...",
#         now published in Phys. Rev. D as "Auxiliary fields approach...") ---
$ws.Range("A15").Value = 45209
$ws.Range("B15").Value = "Auxiliary fields approach to shift-symmetric theories: the φ^4 derivative theory and the crumpled-to-flat transition of membranes at two-loop order"
$ws.Range("C15").Value = "L. Delzescaux, C. Duclut, D. Mouhanna, and M. Tissier"
$ws.Range("D15").Value = "Phys. Rev. D"
$ws.Range("E15").Value = 108
$ws.Range("F15").Value = "L081702"
$ws.Range("G15").Value = "We introduce a technique relying on the use of auxiliary fields in order to eliminate explicit field-derivatives that plague the high orders renormalization group treatment of shift-symmetric (derivative) theories. This technique simplifies drastically the computation of fluctuations in such theories. This is illustrated by deriving the two-loop renormalization group equations—and the three-loop anomalous dimension—of the φ^4 derivative theory in D=4−ε, which is also relevant to describe the crumpled-to-flat transition of polymerized membranes. Some features of this transition are provided."
$ws.Range("H15").Value = "delzescaux2023lagrange"
$ws.Range("I15").Value = "https://link.aps.org/doi/10.1103/PhysRevD.108.L081702"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "2307.00600"
$ws.Range("J15").Style = "Normal"
$ws.Range("K15").Value = "statisticalPhysics, phaseTransition, renormalizationGroup"

# --- 3. Row 20 (the old row 19 that got pushed down by the insert) still carries
#         the stale "preprint" data for delzescaux2023lagrange; overwrite it with
#         the brand-new "duclut2023probe" publication record instead ---
$ws.Range("A20").Value = 45211
$ws.Range("B20").Value = "Probe particles in odd active viscoelastic fluids: how activity and dissipation determine linear stability"
$ws.Range("C20").Value = "C. Duclut, S. Bo, R. Lier, J. Armas, F. Jülicher, P. Surówka"
$ws.Range("D20").Value = "preprint"
$ws.Range("E20").Value = "/"
$ws.Range("F20").Value = "/"
$ws.Range("G20").Value = "Odd viscoelastic materials are constrained by fewer symmetries than their even counterparts. The breaking of these symmetries allow these materials to exhibit different features, which have attracted considerable attention in recent years. Immersing a bead in such complex fluids allows for probing their physical properties, highlighting signatures of their oddity and exploring consequences of these broken symmetries. We present the conditions under which the activity of an odd viscoelastic fluid can give rise to linear instabilities in the motion of the probe particle and unveil how the features of the probe particle dynamics depend on the oddity and activity of the viscoelastic medium in which it is immersed."
$ws.Range("H20").Value = "duclut2023probe"
$ws.Range("I20").Value = "/"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "2310.08640"
$ws.Range("J20").Style = "Normal"
$ws.Range("K20").Value = "oddViscosity, oddElasticity, viscoelastic, thermodynamics, hydrodynamics"

# --- 4. Restore last-selected cell as recorded in the saved view state ---
$ws.Range("C26").Select()
